# regen save_data to use K instead of Strike#, regen std/mean, calc and write s_vals
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New K (strikeouts) values replacing the old Strike# counts in column G,
# rows 2-33 (one row per game).
$kValues = @{
    2  = 8
    3  = 2
    4  = 4
    5  = 5
    6  = 2
    7  = 4
    8  = 6
    9  = 6
    10 = 4
    11 = 6
    12 = 7
    13 = 0
    14 = 4
    15 = 4
    16 = 6
    17 = 6
    18 = 2
    19 = 9
    20 = 10
    21 = 3
    22 = 9
    23 = 6
    24 = 8
    25 = 4
    26 = 8
    27 = 8
    28 = 4
    29 = 7
    30 = 4
    31 = 5
    32 = 5
    33 = 2
}

foreach ($row in $kValues.Keys) {
    $ws.Cells.Item($row, 7).Value = $kValues[$row]
}
